$wb = $excel.ActiveWorkbook

# The handback file (abf76a13-...) has been generated/processed: its status
# moves from "Ready for handoff" to "Handed back: in sync with en-US" on the
# Overview sheet, and on each language sheet its Status + Latest Handback
# DateTime are refreshed and the stale "version mismatch" Error Detail is
# cleared.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("L3").Value = "2017-02-21 03:18:16"
$wsZhCn.Range("R3").Value = ""

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("L3").Value = "2017-02-21 03:18:39"
$wsDeDe.Range("R3").Value = ""
